$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 01:05"

# --- Argentina overtakes Mexico: swap country labels for rows 11/12, keep updated numbers ---
$ws.Range("A11").Value = "Argentina"
$ws.Range("A12").Value = "Mexico"

# --- Uruguay overtakes Yemen: swap country labels for rows 155/156 ---
$ws.Range("A155").Value = "Uruguay"
$ws.Range("A156").Value = "Yemen"

# --- Updated daily statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4
$ws.Range("B4").Value = 7443353
$ws.Range("C4").Value = 37000
$ws.Range("D4").Value = 4687574
$ws.Range("E4").Value = 2544109
$ws.Range("G4").Value = 885
$ws.Range("H4").Value = 211670

# Row 6
$ws.Range("B6").Value = 4810935
$ws.Range("C6").Value = 30618
$ws.Range("D6").Value = 4180376
$ws.Range("E6").Value = 486597
$ws.Range("G6").Value = 952
$ws.Range("H6").Value = 143962

# Row 8
$ws.Range("B8").Value = 829679
$ws.Range("C8").Value = 5637
$ws.Range("D8").Value = 743653
$ws.Range("E8").Value = 60028
$ws.Range("G8").Value = 170
$ws.Range("H8").Value = 25998

# Row 11
$ws.Range("B11").Value = 751001
$ws.Range("C11").Value = 14392
$ws.Range("D11").Value = 594645
$ws.Range("E11").Value = 139419
$ws.Range("G11").Value = 418
$ws.Range("H11").Value = 16937

# Row 12
$ws.Range("B12").Value = 738163
$ws.Range("C12").Value = 4446
$ws.Range("D12").Value = 530945
$ws.Range("E12").Value = 130055
$ws.Range("G12").Value = 560
$ws.Range("H12").Value = 77163

# Row 25
$ws.Range("B25").Value = 292911
$ws.Range("C25").Value = 2445
$ws.Range("E25").Value = 27340

# Row 41
$ws.Range("B41").Value = 103198
$ws.Range("C41").Value = 119
$ws.Range("D41").Value = 96494
$ws.Range("E41").Value = 774
$ws.Range("G41").Value = 16
$ws.Range("H41").Value = 5930

# Row 48
$ws.Range("B48").Value = 83010
$ws.Range("C48").Value = 516
$ws.Range("D48").Value = 76025
$ws.Range("E48").Value = 5421
$ws.Range("G48").Value = 7
$ws.Range("H48").Value = 1564

# Row 56
$ws.Range("B56").Value = 70864
$ws.Range("C56").Value = 442
$ws.Range("D56").Value = 64838
$ws.Range("E56").Value = 5775
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 251

# Row 58
$ws.Range("B58").Value = 58848
$ws.Range("C58").Value = 201
$ws.Range("D58").Value = 50358
$ws.Range("E58").Value = 7378
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1112

# Row 84
$ws.Range("B84").Value = 20833
$ws.Range("C84").Value = 286
$ws.Range("D84").Value = 14634
$ws.Range("E84").Value = 5374
$ws.Range("G84").Value = 12
$ws.Range("H84").Value = 825

# Row 93
$ws.Range("B93").Value = 14027
$ws.Range("C93").Value = 113
$ws.Range("E93").Value = 2563

# Row 109
$ws.Range("B109").Value = 8766
$ws.Range("C109").Value = 14
$ws.Range("D109").Value = 8005
$ws.Range("E109").Value = 707

# Row 113
$ws.Range("B113").Value = 7838
$ws.Range("C113").Value = 1
$ws.Range("D113").Value = 6303
$ws.Range("E113").Value = 1307

# Row 131
$ws.Range("B131").Value = 4829
$ws.Range("C131").Value = 23
$ws.Range("D131").Value = 1914
$ws.Range("E131").Value = 2853

# Row 133
$ws.Range("B133").Value = 4531
$ws.Range("C133").Value = 68
$ws.Range("E133").Value = 1895
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 76

# Row 134
$ws.Range("B134").Value = 4200
$ws.Range("C134").Value = 52
$ws.Range("D134").Value = 1103
$ws.Range("E134").Value = 2897
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 200

# Row 155
$ws.Range("B155").Value = 2046
$ws.Range("C155").Value = 13
$ws.Range("D155").Value = 1791
$ws.Range("E155").Value = 207
$ws.Range("H155").Value = 48

# Row 156
$ws.Range("B156").Value = 2034
$ws.Range("C156").Value = 3
$ws.Range("D156").Value = 1286
$ws.Range("E156").Value = 161
$ws.Range("H156").Value = 587

# Row 169
$ws.Range("D169").Value = 886
$ws.Range("E169").Value = 10

